$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.823.65"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.478.85"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'319.42"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'93.44"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'0.0881"
$ws.Range("E10").Value = "  +11.20%  "
$ws.Range("D11").Value = "'33.37"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "2.860.62"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "'6.94"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'15.75"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "2.483.01"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "41.770.94"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").Value = "'6.47"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'71.24"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").Value = "'240.08"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'2.77"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'24.80"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'9.80"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "'37.15"
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("D31").Value = "'157.47"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").Value = "'5.53"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "'0.0768"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "'17.53"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  +5.01%  "
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "'0.105"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  +6.68%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "2.004.72"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").Value = "'19.15"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").Value = "'9.49"
$ws.Range("E47").Value = "  +4.48%  "
$ws.Range("D48").Value = "2.717.40"
$ws.Range("D49").Value = "'98.54"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("D50").Value = "'75.12"
$ws.Range("E50").Value = "  +4.85%  "
$ws.Range("D51").Value = "'67.59"
$ws.Range("E51").Value = "  +0.53%  "
